$d = $word.ActiveDocument

# Change 1: "By Mo Mo" -> "By Golf"
$d.Content.Find.Execute("By Mo Mo", $false, $false, $false, $false, $false,
                         $true, 1, $false, "By Golf", 2)

# Change 2: "Title 5  Test 123456789" -> same text, but split into separate
# runs by Word's grammar checker (proofErr gramStart/gramEnd around "5  Test").
# The visible text content is unchanged: "Title 5  Test 123456789"
$d.Content.Find.Execute("Title 5  Test 123456789", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Title 5  Test 123456789", 2)
